$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet2" to "Sheet1"
$ws.Name = "Sheet1"

# Update header labels (shared strings)
$ws.Range("B1").Value = "Intervals"
$ws.Range("C1").Value = "Number Of Count"

# Update existing row 2 values
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1

# Add new row 3 with values
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 1

# Copy the style of A2 onto A3 so the new cell matches (border/bold/centered)
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
